# Update the Turkey MSME summary numbers to their more-precise values.
# (Matches upstream "Update country data files" re-export: the previously
# rounded-to-1-decimal figures gain a second decimal, and the "Employment
# (% of total)" MSMEs figure for the SME-Associations table gets its own
# distinct value - 69.17 - instead of re-using the Statistical-Institution
# table's 69.2/69.22 string.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @{
    # Source Type: Statistical Institution table (rows 10-14)
    "B11" = "34.52"   # Enterprises density - Micro
    "C11" = "0.89"    # Enterprises density - SMEs
    "D11" = "35.41"   # Enterprises density - MSMEs
    "B12" = "32.12"   # Employment (% of total) - Micro
    "C12" = "37.11"   # Employment (% of total) - SMEs
    "D12" = "69.22"   # Employment (% of total) - MSMEs
    "B14" = "97.33"   # Enterprises (% of total) - Micro
    "C14" = "2.52"    # Enterprises (% of total) - SMEs
    "D14" = "99.85"   # Enterprises (% of total) - MSMEs

    # Source Type: SME Associations table (rows 32-35)
    "B33" = "32.27"   # Enterprises density - Micro
    "C33" = "1.86"    # Enterprises density - SMEs
    "D33" = "34.13"   # Enterprises density - MSMEs
    "B34" = "38.54"   # Employment (% of total) - Micro
    "C34" = "30.63"   # Employment (% of total) - SMEs
    "D34" = "69.17"   # Employment (% of total) - MSMEs (now its own value)
}

foreach ($addr in $edits.Keys) {
    $cell = $ws.Range($addr)
    # Write the new reading as a text-producing formula, then flatten it to
    # a plain value in place. This keeps the cell's stored type a literal
    # string (matching the source workbook's shared-string cells) instead
    # of having these numeric-looking readings coerced into real numbers,
    # without touching the cell's number format / style.
    $cell.Formula = "=""" + $edits[$addr] + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}
$excel.CutCopyMode = $false
